$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10586
$ws.Range("H62").Value = 3062.5
$ws.Range("I62").Value = 1125
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 1125
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -501
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3062.5
$ws.Range("I65").Value = 1125
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 5625
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -2505
$ws.Range("N65").Value = -31240
$ws.Range("H100").Value = 1676.8334
$ws.Range("I100").Value = 1028.375
$ws.Range("K100").Value = 1028.375
$ws.Range("M100").Value = -487.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2258.4211
$ws.Range("I61").Value = 1750.625
$ws.Range("J61").Value = 4966.6665
$ws.Range("K61").Value = 1750.625
$ws.Range("L61").Value = 4966.6665
$ws.Range("M61").Value = -1538.625
$ws.Range("N61").Value = -5390.6665
$ws.Range("H63").Value = 3320.6
$ws.Range("I63").Value = 2386.5715
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 2386.5715
$ws.Range("L63").Value = 5500
$ws.Range("M63").Value = -1700.5715
$ws.Range("N63").Value = -6872
$ws.Range("H66").Value = 3320.6
$ws.Range("I66").Value = 2386.5715
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 11932.8575
$ws.Range("L66").Value = 27500
$ws.Range("M66").Value = -8500.8575
$ws.Range("N66").Value = -34364
$ws.Range("H132").Value = 46544
$ws.Range("I132").Value = 74107.64
$ws.Range("J132").Value = 3667.2222
$ws.Range("K132").Value = 222322.92
$ws.Range("L132").Value = 11001.6666
$ws.Range("M132").Value = -219792.92
$ws.Range("N132").Value = -16061.6666
$ws.Range("H136").Value = 2258.4211
$ws.Range("I136").Value = 1750.625
$ws.Range("J136").Value = 4966.6665
$ws.Range("K136").Value = 5251.875
$ws.Range("L136").Value = 14899.9995
$ws.Range("M136").Value = -2701.875
$ws.Range("N136").Value = -19999.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5090.48
$ws.Range("I107").Value = 5343
$ws.Range("J107").Value = 2186.5
$ws.Range("K107").Value = 5343
$ws.Range("L107").Value = 2186.5
$ws.Range("M107").Value = -3423
$ws.Range("N107").Value = -6026.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2254.025
$ws.Range("I31").Value = 1437.7142
$ws.Range("J31").Value = 3156.2632
$ws.Range("K31").Value = 1437.7142
$ws.Range("L31").Value = 3156.2632
$ws.Range("M31").Value = -1142.7142
$ws.Range("N31").Value = -3746.2632
$ws.Range("H34").Value = 2254.025
$ws.Range("I34").Value = 1437.7142
$ws.Range("J34").Value = 3156.2632
$ws.Range("K34").Value = 1437.7142
$ws.Range("L34").Value = 3156.2632
$ws.Range("M34").Value = -1235.7142
$ws.Range("N34").Value = -3560.2632
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("K39").Value = 2000
$ws.Range("M39").Value = -1609
$ws.Range("H49").Value = 2000
$ws.Range("I49").Value = 2000
$ws.Range("K49").Value = 2000
$ws.Range("M49").Value = -1818
$ws.Range("H86").Value = 211099.45
$ws.Range("I86").Value = 335112.47
$ws.Range("J86").Value = 4411.1113
$ws.Range("K86").Value = 335112.47
$ws.Range("L86").Value = 4411.1113
$ws.Range("M86").Value = -333989.47
$ws.Range("N86").Value = -6657.1113
$ws.Range("H89").Value = 211099.45
$ws.Range("I89").Value = 335112.47
$ws.Range("J89").Value = 4411.1113
$ws.Range("K89").Value = 1675562.35
$ws.Range("L89").Value = 22055.5565
$ws.Range("M89").Value = -1669946.35
$ws.Range("N89").Value = -33287.5565
$ws.Range("H107").Value = 1241.0625
$ws.Range("I107").Value = 980
$ws.Range("J107").Value = 2024.25
$ws.Range("K107").Value = 980
$ws.Range("L107").Value = 2024.25
$ws.Range("M107").Value = 940
$ws.Range("N107").Value = -5864.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 545
$ws.Range("I13").Value = 110
$ws.Range("J13").Value = 980
$ws.Range("K13").Value = 330
$ws.Range("L13").Value = 2940
$ws.Range("M13").Value = -162
$ws.Range("N13").Value = -3276
$ws.Range("H136").Value = 1130.4286
$ws.Range("I136").Value = 1037.6666
$ws.Range("K136").Value = 3112.9998
$ws.Range("M136").Value = 1987.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1840
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6240
$ws.Range("H122").Value = 2080.2
$ws.Range("I122").Value = 1917.8823
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5753.6469
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3303.6469
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 2264.5881
$ws.Range("I126").Value = 1909.091
$ws.Range("J126").Value = 2916.3333
$ws.Range("K126").Value = 5727.272999999999
$ws.Range("L126").Value = 8748.999899999999
$ws.Range("M126").Value = -3257.272999999999
$ws.Range("N126").Value = -13688.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1504.5333
$ws.Range("I40").Value = 1230.25
$ws.Range("J40").Value = 2601.6667
$ws.Range("K40").Value = 1230.25
$ws.Range("L40").Value = 2601.6667
$ws.Range("M40").Value = -1094.25
$ws.Range("N40").Value = -2873.6667
$ws.Range("H61").Value = 3017.3333
$ws.Range("I61").Value = 2520.8
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 2520.8
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -2318.8
$ws.Range("N61").Value = -5904
$ws.Range("H93").Value = 1018.2778
$ws.Range("I93").Value = 1034.6364
$ws.Range("J93").Value = 992.5714
$ws.Range("K93").Value = 1034.6364
$ws.Range("L93").Value = 992.5714
$ws.Range("M93").Value = 213.3635999999999
$ws.Range("N93").Value = -3488.5714
$ws.Range("H100").Value = 2116
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2145
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2145
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3227
$ws.Range("H113").Value = 3017.3333
$ws.Range("I113").Value = 2520.8
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 2520.8
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -350.8000000000002
$ws.Range("N113").Value = -9840
$ws.Range("H122").Value = 4999.269
$ws.Range("I122").Value = 5175.1
$ws.Range("J122").Value = 4413.1665
$ws.Range("K122").Value = 15525.3
$ws.Range("L122").Value = 13239.4995
$ws.Range("M122").Value = -13075.3
$ws.Range("N122").Value = -18139.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 693.8
$ws.Range("I17").Value = 693.8
$ws.Range("K17").Value = 693.8
$ws.Range("M17").Value = -521.8
